# Scheduled-runner update: refresh currentAveragePrice / Leve price & profit
# columns (H-N) on the Sephirot Profits workbook across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 14042.8
$ws.Range("I38").Value = 12553.5
$ws.Range("K38").Value = 37660.5
$ws.Range("M38").Value = -37288.5

# Row 87
$ws.Range("H87").Value = 97354
$ws.Range("J87").Value = 97354
$ws.Range("L87").Value = 97354
$ws.Range("N87").Value = -99850

# Row 90
$ws.Range("H90").Value = 97354
$ws.Range("J90").Value = 97354
$ws.Range("L90").Value = 292062
$ws.Range("N90").Value = -304542

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 102
$ws.Range("H102").Value = 2582
$ws.Range("I102").Value = 2582
$ws.Range("K102").Value = 2582
$ws.Range("M102").Value = -960

# Row 110
$ws.Range("H110").Value = 943.1111
$ws.Range("J110").Value = 696.3333
$ws.Range("L110").Value = 696.3333
$ws.Range("N110").Value = -4786.3333

# Row 122
$ws.Range("H122").Value = 3130.5
$ws.Range("I122").Value = 2507.3333
$ws.Range("K122").Value = 7521.999899999999
$ws.Range("M122").Value = -5071.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2354.4285
$ws.Range("I86").Value = 2580.1667
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2580.1667
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -1457.1667
$ws.Range("N86").Value = -3246

# Row 89
$ws.Range("H89").Value = 2354.4285
$ws.Range("I89").Value = 2580.1667
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 12900.8335
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -7284.833500000001
$ws.Range("N89").Value = -16232

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 3571.3333
$ws.Range("I105").Value = 2985.111
$ws.Range("K105").Value = 2985.111
$ws.Range("M105").Value = -1238.111

# Row 132
$ws.Range("H132").Value = 3616.25
$ws.Range("I132").Value = 1233.5
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 3700.5
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1170.5
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1119.6666
$ws.Range("J113").Value = 1119.6666
$ws.Range("L113").Value = 3358.9998
$ws.Range("N113").Value = -7698.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 100.6
$ws.Range("I2").Value = 50.75
$ws.Range("K2").Value = 50.75
$ws.Range("M2").Value = 62.25

# Row 97
$ws.Range("H97").Value = 1955
$ws.Range("J97").Value = 400
$ws.Range("L97").Value = 400
$ws.Range("N97").Value = -1392

# Row 132
$ws.Range("H132").Value = 2532.5881
$ws.Range("I132").Value = 1932.5714
$ws.Range("K132").Value = 5797.7142
$ws.Range("M132").Value = -3267.7142

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1244.2727
$ws.Range("I22").Value = 581.3333
$ws.Range("J22").Value = 2039.8
$ws.Range("K22").Value = 581.3333
$ws.Range("L22").Value = 2039.8
$ws.Range("M22").Value = -286.3333
$ws.Range("N22").Value = -2629.8

# Row 27
$ws.Range("H27").Value = 1244.2727
$ws.Range("I27").Value = 581.3333
$ws.Range("J27").Value = 2039.8
$ws.Range("K27").Value = 581.3333
$ws.Range("L27").Value = 2039.8
$ws.Range("M27").Value = -474.3333
$ws.Range("N27").Value = -2253.8

# Row 46
$ws.Range("H46").Value = 2666.5557
$ws.Range("I46").Value = 1999.8572
$ws.Range("K46").Value = 1999.8572
$ws.Range("M46").Value = -1811.8572

# Row 55
$ws.Range("H55").Value = 3353
$ws.Range("I55").Value = 3353
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3353
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = -3180
$ws.Range("M55").ClearContents()

# Row 61
$ws.Range("H61").Value = 679.6
$ws.Range("I61").Value = 550
$ws.Range("K61").Value = 550
$ws.Range("M61").Value = -348

# Row 68
$ws.Range("H68").Value = 1443.5
$ws.Range("I68").Value = 1443.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1443.5
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = -694.5
$ws.Range("M68").ClearContents()

# Row 71
$ws.Range("H71").Value = 1443.5
$ws.Range("I71").Value = 1443.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 7217.5
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = -3473.5
$ws.Range("M71").ClearContents()

# Row 82
$ws.Range("H82").Value = 1425
$ws.Range("I82").Value = 854.1667
$ws.Range("J82").Value = 2566.6667
$ws.Range("K82").Value = 854.1667
$ws.Range("L82").Value = 2566.6667
$ws.Range("M82").Value = -493.1667
$ws.Range("N82").Value = -3288.6667

# Row 85
$ws.Range("H85").Value = 1425
$ws.Range("I85").Value = 854.1667
$ws.Range("J85").Value = 2566.6667
$ws.Range("K85").Value = 854.1667
$ws.Range("L85").Value = 2566.6667
$ws.Range("M85").Value = 393.8333
$ws.Range("N85").Value = -5062.6667

# Row 100
$ws.Range("H100").Value = 1894.6
$ws.Range("I100").Value = 1871.25
$ws.Range("J100").Value = 1988
$ws.Range("K100").Value = 1871.25
$ws.Range("L100").Value = 1988
$ws.Range("M100").Value = -1330.25
$ws.Range("N100").Value = -3070

# Row 113
$ws.Range("H113").Value = 679.6
$ws.Range("I113").Value = 550
$ws.Range("K113").Value = 550
$ws.Range("M113").Value = 1620

# Row 122
$ws.Range("H122").Value = 2997.25
$ws.Range("I122").Value = 2996.3333
$ws.Range("K122").Value = 8988.999899999999
$ws.Range("M122").Value = -6538.999899999999

# Row 132
$ws.Range("H132").Value = 5666
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 874.55554
$ws.Range("I107").Value = 721
$ws.Range("J107").Value = 997.4
$ws.Range("K107").Value = 2163
$ws.Range("L107").Value = 2992.2
$ws.Range("M107").Value = -243
$ws.Range("N107").Value = -6832.2

# Row 136
$ws.Range("H136").Value = 463.33334
$ws.Range("J136").Value = 500
$ws.Range("L136").Value = 1500
$ws.Range("N136").Value = -6600

Write-Host "Applied 161 cell updates across 31 rows."
